$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophe forces the Price (D) values to stay plain text so
# strings like "338.20" or "0.00001093" keep their exact formatting
# instead of being coerced into numbers (losing trailing zeros / going
# scientific). It is stripped by Excel and not stored in the cell.

$ws.Cells.Item(2, 4).Value = "'28.062.14"
$ws.Cells.Item(2, 5).Value = "  -0.01%  "
$ws.Cells.Item(3, 4).Value = "'1.821.91"
$ws.Cells.Item(3, 5).Value = "  +2.76%  "
$ws.Cells.Item(4, 4).Value = "'1.003"
$ws.Cells.Item(4, 5).Value = "  -0.40%  "
$ws.Cells.Item(5, 4).Value = "'338.20"
$ws.Cells.Item(5, 5).Value = "  -0.23%  "
$ws.Cells.Item(6, 4).Value = "'1.002"
$ws.Cells.Item(6, 5).Value = "  -0.11%  "
$ws.Cells.Item(7, 4).Value = "'0.4223"
$ws.Cells.Item(7, 5).Value = "  +10.76%  "
$ws.Cells.Item(8, 4).Value = "'0.3527"
$ws.Cells.Item(8, 5).Value = "  +3.29%  "
$ws.Cells.Item(9, 4).Value = "'45.66"
$ws.Cells.Item(9, 5).Value = "  -2.78%  "
$ws.Cells.Item(10, 4).Value = "'1.154"
$ws.Cells.Item(10, 5).Value = "  +1.22%  "
$ws.Cells.Item(11, 4).Value = "'0.07502"
$ws.Cells.Item(11, 5).Value = "  +1.87%  "
$ws.Cells.Item(12, 4).Value = "'23.03"
$ws.Cells.Item(12, 5).Value = "  -1.52%  "
$ws.Cells.Item(13, 4).Value = "'1.005"
$ws.Cells.Item(13, 5).Value = "  +0.12%  "
$ws.Cells.Item(14, 4).Value = "'6.297"
$ws.Cells.Item(14, 5).Value = "  -1.28%  "
$ws.Cells.Item(15, 4).Value = "'7.329"
$ws.Cells.Item(15, 5).Value = "  +0.55%  "
$ws.Cells.Item(16, 4).Value = "'1.818.23"
$ws.Cells.Item(16, 5).Value = "  +2.23%  "
$ws.Cells.Item(17, 4).Value = "'0.00001093"
$ws.Cells.Item(17, 5).Value = "  +1.73%  "
$ws.Cells.Item(18, 5).Value = "  +0.67%  "
$ws.Cells.Item(19, 4).Value = "'82.59"
$ws.Cells.Item(19, 5).Value = "  +0.38%  "
$ws.Cells.Item(20, 4).Value = "'1.001"
$ws.Cells.Item(20, 5).Value = "  -0.30%  "
$ws.Cells.Item(21, 4).Value = "'17.43"
$ws.Cells.Item(21, 5).Value = "  +0.78%  "
$ws.Cells.Item(22, 5).Value = "  +0.57%  "
$ws.Cells.Item(23, 4).Value = "'28.059.23"
$ws.Cells.Item(23, 5).Value = "  -0.19%  "
$ws.Cells.Item(24, 4).Value = "'11.93"
$ws.Cells.Item(24, 5).Value = "  -1.13%  "
$ws.Cells.Item(25, 4).Value = "'2.400"
$ws.Cells.Item(25, 5).Value = "  +0.33%  "
$ws.Cells.Item(26, 4).Value = "'2.491"
$ws.Cells.Item(26, 5).Value = "  +4.14%  "
$ws.Cells.Item(27, 4).Value = "'20.86"
$ws.Cells.Item(27, 5).Value = "  +1.08%  "
$ws.Cells.Item(28, 4).Value = "'156.70"
$ws.Cells.Item(28, 5).Value = "  +2.10%  "
$ws.Cells.Item(29, 4).Value = "'2.027.80"
$ws.Cells.Item(29, 5).Value = "  +2.50%  "
$ws.Cells.Item(30, 5).Value = "  -7.18%  "
$ws.Cells.Item(31, 4).Value = "'133.51"
$ws.Cells.Item(31, 5).Value = "  -0.88%  "
$ws.Cells.Item(32, 4).Value = "'4.086"
$ws.Cells.Item(32, 5).Value = "  +1.64%  "
$ws.Cells.Item(33, 4).Value = "'6.022"
$ws.Cells.Item(33, 5).Value = "  -0.27%  "
$ws.Cells.Item(34, 4).Value = "'0.09174"
$ws.Cells.Item(34, 5).Value = "  +3.52%  "
$ws.Cells.Item(35, 4).Value = "'12.41"
$ws.Cells.Item(35, 5).Value = "  -1.98%  "
$ws.Cells.Item(36, 4).Value = "'0.02363"
$ws.Cells.Item(36, 5).Value = "  -1.42%  "
$ws.Cells.Item(37, 4).Value = "'0.06326"
$ws.Cells.Item(37, 5).Value = "  -0.10%  "
$ws.Cells.Item(38, 4).Value = "'0.6680"
$ws.Cells.Item(38, 5).Value = "  -1.64%  "
$ws.Cells.Item(39, 4).Value = "'5.263"
$ws.Cells.Item(39, 5).Value = "  -0.78%  "
$ws.Cells.Item(40, 4).Value = "'0.2183"
$ws.Cells.Item(40, 5).Value = "  +1.68%  "
$ws.Cells.Item(41, 4).Value = "'1.513"
$ws.Cells.Item(41, 5).Value = "  +1.05%  "
$ws.Cells.Item(42, 4).Value = "'1.221"
$ws.Cells.Item(42, 5).Value = "  -1.35%  "
$ws.Cells.Item(43, 4).Value = "'8.135"
$ws.Cells.Item(43, 5).Value = "  -1.13%  "
$ws.Cells.Item(44, 4).Value = "'14.29"
$ws.Cells.Item(44, 5).Value = "  +1.18%  "
$ws.Cells.Item(45, 4).Value = "'1.001"
$ws.Cells.Item(45, 5).Value = "  -0.16%  "
$ws.Cells.Item(46, 4).Value = "'0.6169"
$ws.Cells.Item(46, 5).Value = "  -0.89%  "
$ws.Cells.Item(47, 4).Value = "'3.880"
$ws.Cells.Item(47, 5).Value = "  +0.55%  "
$ws.Cells.Item(48, 4).Value = "'128.79"
$ws.Cells.Item(48, 5).Value = "  -2.64%  "
$ws.Cells.Item(49, 4).Value = "'2.068"
$ws.Cells.Item(49, 5).Value = "  +0.63%  "
$ws.Cells.Item(50, 4).Value = "'1.185"
$ws.Cells.Item(50, 5).Value = "  -0.79%  "
$ws.Cells.Item(51, 5).Value = "  -4.77%  "
